$wb = $excel.ActiveWorkbook

# This script adds a new handback entry for file
# "cad4860f-87f8-48a8-8617-6517a6f51269.md" to all three report sheets
# (Overview, zh-cn, de-de), mirroring the existing rows for the other
# two tracked files, and appends the corresponding new table row + hyperlinks.

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$ro = $rowOverview.Range.Row

$wsOverview.Cells.Item($ro, 1).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.md"
$wsOverview.Cells.Item($ro, 3).Value = ".md"
$wsOverview.Cells.Item($ro, 5).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item($ro, 6).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item($ro, 7).Value = "2016-08-22 00:55:50"
$wsOverview.Cells.Item($ro, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($ro, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/cad4860f-87f8-48a8-8617-6517a6f51269.md", "", "", "e2e\cad4860f-87f8-48a8-8617-6517a6f51269.md")

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rz = $rowZhCn.Range.Row

$wsZhCn.Cells.Item($rz, 1).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.md"
$wsZhCn.Cells.Item($rz, 2).Value = ".md"
$wsZhCn.Cells.Item($rz, 3).Value = "Handed back: in sync with en-US"
$wsZhCn.Cells.Item($rz, 4).Value = "e2e"
$wsZhCn.Cells.Item($rz, 5).Value = "ht"

$c = $wsZhCn.Cells.Item($rz, 6); $c.Value = "'True"; $c.Style = "Normal"

$wsZhCn.Cells.Item($rz, 7).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.a97fcdc272db142a51d7759af214123c3c45cc1a.zh-cn.xlf"
$wsZhCn.Cells.Item($rz, 8).Value = "2016-08-22 00:55:45"
$wsZhCn.Cells.Item($rz, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($rz, 9).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.md"
$wsZhCn.Cells.Item($rz, 10).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.a97fcdc272db142a51d7759af214123c3c45cc1a.zh-cn.xlf"
$wsZhCn.Cells.Item($rz, 11).Value = "2016-08-22 00:56:06"
$wsZhCn.Cells.Item($rz, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$c = $wsZhCn.Cells.Item($rz, 12); $c.Value = "'"; $c.Style = "Normal"
$c = $wsZhCn.Cells.Item($rz, 13); $c.Value = "'True"; $c.Style = "Normal"
$c = $wsZhCn.Cells.Item($rz, 14); $c.Value = "'"; $c.Style = "Normal"
$c = $wsZhCn.Cells.Item($rz, 15); $c.Value = "'False"; $c.Style = "Normal"
$c = $wsZhCn.Cells.Item($rz, 16); $c.Value = "'"; $c.Style = "Normal"

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item($rz, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/cad4860f-87f8-48a8-8617-6517a6f51269.md", "", "", "cad4860f-87f8-48a8-8617-6517a6f51269.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item($rz, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/cad4860f-87f8-48a8-8617-6517a6f51269.md", "", "", "cad4860f-87f8-48a8-8617-6517a6f51269.md")

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rd = $rowDeDe.Range.Row

$wsDeDe.Cells.Item($rd, 1).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.md"
$wsDeDe.Cells.Item($rd, 2).Value = ".md"
$wsDeDe.Cells.Item($rd, 3).Value = "Handed back: in sync with en-US"
$wsDeDe.Cells.Item($rd, 4).Value = "e2e"
$wsDeDe.Cells.Item($rd, 5).Value = "ht"

$c = $wsDeDe.Cells.Item($rd, 6); $c.Value = "'True"; $c.Style = "Normal"

$wsDeDe.Cells.Item($rd, 7).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.a97fcdc272db142a51d7759af214123c3c45cc1a.de-de.xlf"
$wsDeDe.Cells.Item($rd, 8).Value = "2016-08-22 00:55:50"
$wsDeDe.Cells.Item($rd, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($rd, 9).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.md"
$wsDeDe.Cells.Item($rd, 10).Value = "cad4860f-87f8-48a8-8617-6517a6f51269.a97fcdc272db142a51d7759af214123c3c45cc1a.de-de.xlf"
$wsDeDe.Cells.Item($rd, 11).Value = "2016-08-22 00:56:13"
$wsDeDe.Cells.Item($rd, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$c = $wsDeDe.Cells.Item($rd, 12); $c.Value = "'"; $c.Style = "Normal"
$c = $wsDeDe.Cells.Item($rd, 13); $c.Value = "'True"; $c.Style = "Normal"
$c = $wsDeDe.Cells.Item($rd, 14); $c.Value = "'"; $c.Style = "Normal"
$c = $wsDeDe.Cells.Item($rd, 15); $c.Value = "'False"; $c.Style = "Normal"
$c = $wsDeDe.Cells.Item($rd, 16); $c.Value = "'"; $c.Style = "Normal"

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item($rd, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/cad4860f-87f8-48a8-8617-6517a6f51269.md", "", "", "cad4860f-87f8-48a8-8617-6517a6f51269.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item($rd, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/cad4860f-87f8-48a8-8617-6517a6f51269.md", "", "", "cad4860f-87f8-48a8-8617-6517a6f51269.md")

Write-Output "Added handback row for cad4860f-87f8-48a8-8617-6517a6f51269.md to Overview, zh-cn and de-de sheets"
